$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two sample data rows (rows 2 and 3), keeping only the header row.
$ws.Range("2:3").EntireRow.Delete() | Out-Null

# Insert a new column for "Subject Code" right after "Subject" (before old column E).
# The inserted column automatically inherits the formatting (style) of the column to its left.
$ws.Columns.Item(5).EntireColumn.Insert() | Out-Null
$ws.Range("E1").Value = "Subject Code"
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Insert a new column for "Subject PRE-Requistie Code" right after "Subject PRE-Requistie"
# (before old "Maximum Students", now shifted to column G).
$ws.Columns.Item(7).EntireColumn.Insert() | Out-Null
$ws.Range("G1").Value = "Subject PRE-Requistie Code"

# Update the saved selection to an empty cell below the data, like the source workbook.
$ws.Range("B11").Select() | Out-Null
